$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - HDKSOE / 009540.KS
$ws.Range("D2").Value = 425500
$ws.Range("E2").Value = 46
$ws.Range("F2").Value = 3.78
$ws.Range("N2").Value = 54.86376272656823

# Row 3 - HD HYUNDAI MIPO / 010620.KS
$ws.Range("N3").Value = 54.86376272656823

# Row 4 - Hanwha Ocean / 042660.KS
$ws.Range("D4").Value = 106800
$ws.Range("E4").Value = 17.3
$ws.Range("F4").Value = -0.93
$ws.Range("I4").Value = 83
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 52.7
$ws.Range("N4").Value = 54.86376272656823

# Row 5 - SamsungHvyInd / 010140.KS
$ws.Range("D5").Value = 24800
$ws.Range("E5").Value = 38.5
$ws.Range("F5").Value = 0.8100000000000001
$ws.Range("I5").Value = 60
$ws.Range("K5").Value = 46.5
$ws.Range("N5").Value = 54.86376272656823
